# Updates the cryptocurrency price/volume table (cryptos list refresh).
# Mirrors: 'Updated cryptos list ... with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.918.47"
$ws.Range("E2").Value = "  +5.73%  "
$ws.Range("D3").Value = "2.233.61"
$ws.Range("E3").Value = "  +2.89%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'231.83"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").Value = "'61.37"
$ws.Range("E7").Value = "  -2.68%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "'0.402"
$ws.Range("E9").Value = "  +2.90%  "
$ws.Range("D10").Value = "'59.11"
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("D11").Value = "'0.0892"
$ws.Range("E11").Value = "  +4.10%  "
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "2.566.43"
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("D14").Value = "'15.67"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Value = "'21.80"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "'0.802"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "'5.58"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").Value = "2.252.71"
$ws.Range("E18").Value = "  +3.81%  "
$ws.Range("D19").Value = "41.875.06"
$ws.Range("E19").Value = "  +5.77%  "
$ws.Range("D20").Value = "'72.35"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "0.0₃0891"
$ws.Range("E21").Value = "  -5.16%  "
$ws.Range("D22").Value = "'6.04"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").Value = "'251.31"
$ws.Range("E23").Value = "  +10.30%  "
$ws.Range("D25").Value = "'2.40"
$ws.Range("E25").Value = "  +2.10%  "
$ws.Range("D26").Value = "'2.36"
$ws.Range("E26").Value = "  +2.27%  "
$ws.Range("D27").Value = "'9.71"
$ws.Range("E27").Value = "  +2.38%  "
$ws.Range("E28").Value = "  +3.98%  "
$ws.Range("D29").Value = "'167.24"
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("D30").Value = "'19.98"
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("E31").Value = "  -2.71%  "
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("D33").Value = "'0.122"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").Value = "'4.98"
$ws.Range("E34").Value = "  +5.84%  "
$ws.Range("D35").Value = "'4.65"
$ws.Range("E35").Value = "  +3.15%  "
$ws.Range("D36").Value = "'0.0633"
$ws.Range("E36").Value = "  +3.08%  "
$ws.Range("D37").Value = "'6.66"
$ws.Range("E37").Value = "  -4.53%  "
$ws.Range("D38").Value = "'3.67"
$ws.Range("E38").Value = "  -3.13%  "
$ws.Range("D39").Value = "'2.37"
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("D40").Value = "'0.000254"
$ws.Range("E40").Value = "  +30.41%  "
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("E42").Value = "  +5.76%  "
$ws.Range("D43").Value = "'4.82"
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("D44").Value = "'8.57"
$ws.Range("E44").Value = "  +8.73%  "
$ws.Range("D45").Value = "'0.0983"
$ws.Range("E45").Value = "  +7.32%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'1.22"
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'99.03"
$ws.Range("E47").Value = "  -3.38%  "
$ws.Range("D48").Value = "1.472.74"
$ws.Range("E48").Value = "  -2.81%  "
$ws.Range("D49").Value = "'16.53"
$ws.Range("E49").Value = "  -6.76%  "
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'52.30"
$ws.Range("E51").Value = "  +8.40%  "
